# Update cryptos list (price / 1h volume change) for the scheduled GitHub Actions refresh.
# Note: some "Price" column strings look numeric (e.g. "208.83") but must stay as plain
# text to match the original inlineStr cell type; a leading apostrophe forces Excel's
# COM layer to keep them as text instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.095.75"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "1.568.19"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.64%  "

$ws.Range("D5").Value = "'208.83"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("D8").Value = "'22.05"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").Value = "'0.0597"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "1.564.46"
$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").Value = "'0.520"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "27.093.35"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "'62.00"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "0.0₃0705"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'7.42"
$ws.Range("E18").Value = "  +2.16%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'215.88"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").Value = "'4.15"
$ws.Range("E21").Value = "  +2.27%  "

$ws.Range("D22").Value = "'9.20"
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "'1.94"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'154.27"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "'6.62"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "'15.07"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").Value = "'0.105"
$ws.Range("E27").Value = "  +1.21%  "

$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("E29").Value = "  +4.55%  "

$ws.Range("D30").Value = "'0.0473"
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("D32").Value = "'3.19"
$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("D33").Value = "1.428.85"
$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("E34").Value = "  +12.88%  "

$ws.Range("E35").Value = "  +1.13%  "

$ws.Range("E36").Value = "  +3.55%  "

$ws.Range("D37").Value = "'0.0167"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").Value = "'0.533"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").Value = "'5.84"
$ws.Range("E39").Value = "  +2.67%  "

$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").Value = "'2.38"
$ws.Range("E41").Value = "  +4.65%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("D44").Value = "'64.64"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").Value = "1.706.57"

$ws.Range("D47").Value = "'86.60"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").Value = "'0.0518"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").Value = "'0.0962"

$ws.Range("E51").Value = "  +0.72%  "
